$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert a new record at row 427 (pushing the existing
# rows 427:456 down to 428:457) and populate it with this week's reading.
$ws.Rows(427).Insert()

$ws.Range("A427").Value = 5
$ws.Range("B427").Value = "Macroferia Regional de Talca"
$ws.Range("C427").Value = "Maule"
$ws.Range("D427").Value = 44585
$ws.Range("E427").Value = 7
$ws.Range("F427").Value = 100112004
$ws.Range("G427").Value = "Cebolla"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "1a nueva(o)"
$ws.Range("J427").Value = 50000
$ws.Range("K427").Value = 1200
$ws.Range("L427").Value = 1200
$ws.Range("M427").Value = 1200
$ws.Range("N427").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O427").Value = "Región del Maule"
$ws.Range("P427").Value = 120
$ws.Range("Q427").Value = 10
$ws.Range("R427").Value = "Hortaliza"
